$wb = $excel.ActiveWorkbook

# --- YDS sheet: append Week 16 play-by-play yardage logs (Run/Pass, OFF/DEF) ---
$ws = $wb.Worksheets.Item("YDS")
$ws.Range("B2").Value = "41 0 14 2 -3 -2 1 3 5 10 3 3 0 9 3 1 6 4 4 6 1 3 8 5 4 3 2 2 3 7 8 6 11 14 5 0 -1 0 6 6 5 11 7 2 0 -3 6 1 1 4 4 -1 1 10 6 4 16 1 7 0 -2 1 6 4 2 1 5 4 1 9 4 1 4 4 10 1 4 -1 5 1 11 9 1 4 10 2 3 2 10 1 4 3 3 24 4 -1 5 7 2 23 8 -3 4 -1 4 6 1 1 3 0 5 11 14 0 0 6 8 9 3 3 9 7 0 4 8 8 1 0 0 6 2 0 7 1 6 6 5 0 1 0 1 1 9 2 0 0 1 25 5 4 1 3 7 3 4 2 13 4 2 2 11 3 6 1 2 14 3 16 0 29 1 1 6 4 0 6 13 5 12 6 0 6 0 8 8 1 2 -1 4 -1 -1 1 1 8 2 6 16 -1 4 22 9 -2 1 3 6 9 6 3 2 3 4 16 3 9 3 16 2 0 15 -3 -4 4 -4 4 6 3 0 14 8 4 0 0 2 3 8 21 2 8 1 6 2 8 1 0 6 2 2 3 2 10 3 1 12 5 2 15 3 8 8 9 2 11 5 -6 3 10 7 5 6 4 5 5 -1 6 11 2 0 7 2 0 0 0 1 2 4 8 4 0 2 6 0 0 3 3 5 4 5 12 -1 9 1 0 1 0 4 9 2 8 7 3 1 13 0 24 -1 7 -2 1 2 2 3 3 8 1 0 5 5 4 6 5 8 2 1 5 1 20 1 1 3 8 3 8 7 0 20 4 5 2 2 3 11 17 1 9 7 13 -1 1 9 3 4 3 3 -2 16 -1 8 0 4 4 11 -1 2 0 1 -7 5 -2 23 3 4 11 2 2 -2 8 3 2 -2 0 0 1 -2 1 12 8 3 -4 4 4 -1 9 2 20 3 0 -1 1"
$ws.Range("C2").Value = "6 1 2 7 0 -2 0 5 1 5 8 9 5 2 2 15 3 3 8 2 3 2 3 1 8 -1 -3 4 9 5 1 5 10 3 3 5 2 4 3 2 12 1 4 2 2 10 12 4 1 3 13 -1 4 1 11 6 0 3 0 12 1 4 7 29 14 3 5 3 19 1 6 2 -6 4 6 5 4 2 5 5 8 4 5 3 -1 1 3 11 2 -2 2 4 9 2 -1 2 4 6 6 3 1 3 -4 0 3 11 3 8 6 5 0 5 -1 5 5 4 2 7 1 1 0 1 1 8 16 6 0 3 3 0 5 8 2 6 36 1 -6 28 1 5 1 5 1 0 7 6 8 3 9 3 15 2 4 5 12 3 9 4 11 6 8 12 2 -1 1 12 19 12 9 5 3 1 2 -3 7 1 1 27 1 -1 3 0 4 2 3 11 5 27 1 5 39 -1 12 -2 0 1 1 3 12 10 2 2 7 2 0 5 8 4 5 4 -2 11 0 11 2 0 0 3 1 2 1 13 -1 3 1 5 10 1 -1 -3 4 1 4 2 5 4 2 19 0 3 8 2 3 4 5 2 5 4 6 5 0 2 1 6 -1 2 2 3 4 -4 0 0 2 4 0 3 7 1 8 4 2 4 2 6 11 19 3 4 1 8 5 0 9 8 21 2 7 1 2 -4 5 10 7 2 14 -1 -6 10 2 5 4 9 0 -6 9 3 3 6 3 0 -4 2 15 1 2 3 6 11 0 11 5 8 11 1 3 14 0 5 3 6 6 2 1 9 7 1 4 7 10 2 3 10 9 2 0 -1 8 0 5 19 11 -2 6 2 6 0 3 1 3 3 5 3 2 2 4 2 1 3 3 -1 25 5 0 2 3 2 6 5 7 3 10 6 2 32 3 7 32 3 3 7 2"
$ws.Range("B3").Value = "4 9 5 10 4 19 -2 10 9 7 3 10 9 8 6 7 10 8 3 10 1 11 5 11 9 9 10 5 6 10 17 -5 3 3 11 2 36 10 5 11 10 13 21 5 0 11 9 9 18 11 10 4 4 21 64 8 8 28 32 6 27 7 4 18 7 2 2 15 8 12 8 9 7 13 14 15 11 20 3 5 3 12 8 11 20 21 8 10 5 15 7 6 5 5 5 3 6 3 17 4 12 4 10 22 20 6 22 12 -6 11 8 5 14 6 11 11 16 6 8 10 4 2 -6 19 22 6 2 9 10 20 5 3 6 25 8 21 17 50 10 11 20 28 5 9 39 16 11 12 6 29 22 22 23 60 10 2 3 12 9 4 23 3 11 13 8 7 12 33 2 17 52 17 16 8 9 0 10 19 11 16 12 7 13 12 7 12 3 5 6 3 10 7 6 10 34 11 2 1 7 12 11 6 6 2 7 19 0 11 19 8 3 6 46 6 54 6 2 5 6 13 5 11 4 3 22 5 12 27 6 8 9 0 11 6 2 9 11 8 12 11 23 11 6 13 3 16 10 10 17 19 4 1 10 0 8 15 9 15 7 5 16 2 15 16 10 16 14 11 5 30 10 14 15"
$ws.Range("C3").Value = "7 67 1 18 4 17 17 1 6 19 56 6 37 3 6 17 15 9 2 2 7 12 10 3 12 1 2 14 8 7 14 14 22 5 10 15 42 7 13 13 21 7 0 3 18 17 7 18 17 6 4 9 13 26 16 23 15 7 21 5 6 33 17 2 20 11 22 5 9 4 11 9 24 13 25 17 3 8 9 9 3 12 8 29 9 29 4 4 9 6 17 5 6 18 3 2 4 6 3 12 7 10 13 4 0 6 32 1 12 5 11 9 4 10 12 14 13 41 8 -4 6 4 21 4 11 9 10 46 2 9 14 0 9 8 9 7 12 0 34 21 -4 19 11 6 23 16 6 15 50 3 83 15 7 11 19 23 2 26 9 11 4 0 7 3 1 3 42 10 9 4 10 2 11 4 12 22 13 11 1 16 5 9 8 -2 11 11 2 3 12 9 15 2 9 8 -2 3 9 4 12 7 21 6 29 0 3 6 16 39 9 10 5 4 -6 9 2 0 -1 17 8 0 11 17 15 7 10 12 23 8 7 12 13 0 4 14 6 6 32 0 7 12 25 12 20 6 13 6 18 38 12 16 8 23 16 12 6 6 8 11 6 4 7 3 5 3 0 2 18 12 4 9 7 7 7 13 7 41 0 9 8 9 13 13 26 3 3 24 9 -3 4 15"

# --- OFF sheet: Week 16 down/distance + play totals ---
$ws = $wb.Worksheets.Item("OFF")
$ws.Range("C2").Value = 194
$ws.Range("D2").Value = 15
$ws.Range("E2").Value = 15
$ws.Range("F2").Value = 72
$ws.Range("G2").Value = 41
$ws.Range("J2").Value = 29
$ws.Range("N2").Value = 29
$ws.Range("O2").Value = 27
$ws.Range("P2").Value = 12
$ws.Range("B3").Value = 10
$ws.Range("C3").Value = 150
$ws.Range("E3").Value = 39
$ws.Range("F3").Value = 74
$ws.Range("G3").Value = 34
$ws.Range("H3").Value = 36
$ws.Range("I3").Value = 49
$ws.Range("J3").Value = 43
$ws.Range("L3").Value = 248
$ws.Range("M3").Value = 152
$ws.Range("Q3").Value = 495

# --- DEF sheet: Week 16 down/distance + play totals ---
$ws = $wb.Worksheets.Item("DEF")
$ws.Range("C2").Value = 188
$ws.Range("D2").Value = 12
$ws.Range("E2").Value = 16
$ws.Range("F2").Value = 58
$ws.Range("G2").Value = 54
$ws.Range("N2").Value = 20
$ws.Range("O2").Value = 24
$ws.Range("C3").Value = 142
$ws.Range("E3").Value = 25
$ws.Range("F3").Value = 83
$ws.Range("G3").Value = 33
$ws.Range("H3").Value = 29
$ws.Range("I3").Value = 41
$ws.Range("J3").Value = 53
$ws.Range("L3").Value = 246
$ws.Range("M3").Value = 168
$ws.Range("Q3").Value = 478

# --- ST sheet: Week 16 special-teams totals + appended KO/PT logs ---
$ws = $wb.Worksheets.Item("ST")
$ws.Range("B2").Value = 58
$ws.Range("D2").Value = 54
$ws.Range("F2").Value = 229
$ws.Range("G2").Value = 221
$ws.Range("H2").Value = 3
$ws.Range("I2").Value = 2
$ws.Range("J2").Value = 113
$ws.Range("K2").Value = 106
$ws.Range("B3").Value = 29
$ws.Range("B4").Value = "70 66 62 67 66 66 62 62 66 64 65 61 56 59 69 62 67 66 66 60 63 64 60 61 66 62 66 46 54"
$ws.Range("B5").Value = "27 18 24 23 21 26 15 15 20 23 17 29 21 19 25 25 35 27 22 21 32 18 19 18 33 23 40 13 16"
$ws.Range("B6").Value = "50 24 13 19 10 22 39 23 31 33 24 32 13 26 27 30 27 31 10 34 24 30 27 25 16 19 26 32 17 17 0 25 19 15 22 27 40 0 42 22 14 18 34 27 22 30 20 8 14"
$ws.Range("D3").Value = "36 42 40 47 40 61 61 52 50 63 48 62 55 44 60 43 53 52 44 42 48 43 43 27 50 28 38 44 41 44 43 49 54 31 53 40 34 50 28 46 51 39 53 38 40 49 42 54 44 45 72 64 28 43"
$ws.Range("D4").Value = "0 0 0 15 0 14 13 11 23 9 14 24 7 0 22 0 0 0 0 9 16 0 0 0 43 0 0 0 21 7 14 6 3 0 7 12 0 0 0 9 0 0 18 0 0 13 0 15 11 3 2 0 0 15"
$ws.Range("D5").Value = "0 0 34 0 97 0 6 -1 0 0 0 0 4 28 0 0"

# --- TURNS sheet: Week 16 turnovers ---
$ws = $wb.Worksheets.Item("TURNS")
$ws.Range("D3").Value = 10
$ws.Range("E3").Value = 8

# --- PEN sheet: Week 16 penalties ---
$ws = $wb.Worksheets.Item("PEN")
$ws.Range("B2").Value = 19
$ws.Range("B3").Value = 14
